# Add three new "Third Party" product rows to the product sheet,
# matching the style/formatting already used by row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Third Party PD - 1"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 12

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "as"
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 12

# Row 15
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Third Party PD - 5"
$ws.Range("C15").Value = 16
$ws.Range("D15").Value = 11
